$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.149639964103699
$ws.Range("B1").Value = 2.906500101089478
$ws.Range("C1").Value = 3.805113315582275
$ws.Range("D1").Value = 3.520087242126465
$ws.Range("E1").Value = 1.203027367591858
